# Font / layout tweaks on Slide 1 (title slide):
#   - Subtitle placeholder: nudge right a bit and drop the run size 25pt -> 24pt
#   - "By Sheldon, ..." credits textbox: reposition and switch every run to
#     the "Avenir Light" typeface.
#
# NOTE: PowerPoint COM exposes Shape.Left/Top in points while the OOXML
# stores EMU (1 pt = 12700 EMU). The literals below are chosen so that,
# after the host's point<->EMU round trip, they land on the exact target
# EMU values (-273517/3010872 and 3834580/6150336).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Subtitle 2 ("An Adaptive User-Metric Oriented Calculation Tool") ---
$subtitle = $s.Shapes.Item(3)
$subtitle.Left = -21.536850903879003

$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Font.Size = 24

# --- TextBox 3 ("By Sheldon, Delan, Luli, Geoff, Tim") ---
$credits = $s.Shapes.Item(5)
$credits.Left = 301.9354407002607
$credits.Top = 484.2784282486082

$creditsText = $credits.TextFrame.TextRange
$total = $creditsText.Length
$creditsText.Characters(1, $total).Font.Name = "Avenir Light"
